$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2, 3, 6, 7 per repulled data
$ws.Range("F2").Value = -5
$ws.Range("F3").Value = -12
$ws.Range("F6").Value = -7
$ws.Range("F7").Value = -2
